# Auto-generated edit script: update cryptos price list
# Applies the per-cell text changes described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'70.162.25"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -0.55%  '
$ws.Range('D3').Value = "'3.509.35"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -0.63%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').Value = "'608.39"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.45%  '
$ws.Range('D6').Value = "'172.79"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.79%  '
$ws.Range('B7').Value = 'XRP'
$ws.Range('C7').Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range('D7').Value = "'0.607"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -1.87%  '
$ws.Range('B8').Value = 'LidoStakedEther'
$ws.Range('C8').Value = 'https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth'
$ws.Range('D8').Value = "'3.500.89"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.76%  '
$ws.Range('D9').Value = "'1.00"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.01%  '
$ws.Range('E10').Value = '  -1.67%  '
$ws.Range('D11').Value = "'7.12"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +6.02%  '
$ws.Range('E12').Value = '  +0.05%  '
$ws.Range('D13').Value = "'46.29"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -1.19%  '
$ws.Range('D14').Value = "'0.0000275"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.93%  '
$ws.Range('D15').Value = "'4.073.87"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.46%  '
$ws.Range('D16').Value = "'8.42"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.63%  '
$ws.Range('D17').Value = "'612.52"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -1.26%  '
$ws.Range('D18').Value = "'3.513.28"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.17%  '
$ws.Range('D19').Value = "'70.175.87"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.59%  '
$ws.Range('E20').Value = '  +2.35%  '
$ws.Range('D21').Value = "'17.77"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +2.11%  '
$ws.Range('D22').Value = "'0.881"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.03%  '
$ws.Range('D23').Value = "'8.99"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -4.17%  '
$ws.Range('D24').Value = "'98.62"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +2.16%  '
$ws.Range('D25').Value = "'15.56"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -1.06%  '
$ws.Range('D26').Value = "'3.75"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -1.22%  '
$ws.Range('D27').Value = "'0.999"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.08%  '
$ws.Range('E28').Value = '  -1.08%  '
$ws.Range('D29').Value = "'33.62"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.68%  '
$ws.Range('D30').Value = "'9.13"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.59%  '
$ws.Range('D31').Value = "'8.09"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -4.35%  '
$ws.Range('D32').Value = "'2.98"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -2.33%  '
$ws.Range('D33').Value = "'1.30"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.57%  '
$ws.Range('E34').Value = '  -2.43%  '
$ws.Range('D35').Value = "'621.88"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +8.66%  '
$ws.Range('D36').Value = "'0.0493"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +5.64%  '
$ws.Range('E37').Value = '  -1.02%  '
$ws.Range('D38').Value = "'10.78"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.15%  '
$ws.Range('D39').Value = "'3.44"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -3.76%  '
$ws.Range('B40').Value = 'OKB'
$ws.Range('C40').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D40').Value = "'56.83"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -1.05%  '
$ws.Range('B41').Value = 'FirstDigitalUSD'
$ws.Range('C41').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D41').Value = "'1.00"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.14%  '
$ws.Range('B42').Value = 'Kaspa'
$ws.Range('C42').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D42').Value = "'0.145"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +1.88%  '
$ws.Range('B43').Value = 'PEPE'
$ws.Range('C43').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D43').Value = "'0.0₃0740"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +4.78%  '
$ws.Range('B44').Value = 'Maker'
$ws.Range('C44').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D44').Value = "'3.359.82"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.36%  '
$ws.Range('D45').Value = "'0.310"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -4.70%  '
$ws.Range('B46').Value = 'ThetaToken'
$ws.Range('C46').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D46').Value = "'2.90"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -4.41%  '
$ws.Range('B47').Value = 'InjectiveProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D47').Value = "'32.10"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -2.76%  '
$ws.Range('D48').Value = "'2.55"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -3.16%  '
$ws.Range('D49').Value = "'0.130"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.38%  '
$ws.Range('D50').Value = "'133.20"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.52%  '
